$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 37 with the new "Head Pat" card
$ws.Range("A37").Value = "Head Pat"
$ws.Range("B37").Value = 0
$ws.Range("C37").Value = "Common"
$ws.Range("D37").Value = "Attack"
$ws.Range("F37").Value = "Deal !D! damage. Apply !M! Passivity."
$ws.Range("E37").Value = "Passivity"

# Update the view selection to match the final saved state
$ws.Range("G38:G39").Select()
